$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / naturally-text values ---
$ws.Range("D2").Value = '66.984.35'
$ws.Range("E2").Value = '  +0.13%  '

$ws.Range("D3").Value = '3.476.14'
$ws.Range("E3").Value = '  +0.73%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("E5").Value = '  -0.72%  '

$ws.Range("E6").Value = '  -0.54%  '

$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("E8").Value = '  -1.24%  '

$ws.Range("D9").Value = '3.473.65'
$ws.Range("E9").Value = '  +0.69%  '

$ws.Range("E10").Value = '  -2.17%  '

$ws.Range("E11").Value = '  -0.28%  '

$ws.Range("E12").Value = '  -1.77%  '

$ws.Range("D13").Value = '4.073.23'
$ws.Range("E13").Value = '  +0.43%  '

$ws.Range("E14").Value = '  -3.02%  '

$ws.Range("E15").Value = '  -0.86%  '

$ws.Range("D16").Value = '66.938.36'
$ws.Range("E16").Value = '  +0.04%  '

$ws.Range("E17").Value = '  -1.07%  '

$ws.Range("D18").Value = '3.468.17'
$ws.Range("E18").Value = '  +0.28%  '

$ws.Range("E19").Value = '  -2.91%  '

$ws.Range("E20").Value = '  -1.14%  '

$ws.Range("E21").Value = '  -1.02%  '

$ws.Range("E22").Value = '  -0.39%  '

$ws.Range("E23").Value = '  +2.51%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("E24").Value = '  +0.74%  '

$ws.Range("B25").Value = 'LEO'
$ws.Range("C25").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("E25").Value = '  +0.58%  '

$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("E26").Value = '  +0.23%  '

$ws.Range("E27").Value = '  +0.32%  '

$ws.Range("E28").Value = '  -3.49%  '

$ws.Range("E29").Value = '  +1.19%  '

$ws.Range("E30").Value = '  -0.19%  '

$ws.Range("E31").Value = '  +4.25%  '

$ws.Range("E32").Value = '  -3.70%  '

$ws.Range("E33").Value = '  -2.00%  '

$ws.Range("E34").Value = '  -3.93%  '

$ws.Range("E36").Value = '  -1.49%  '

$ws.Range("E37").Value = '  +0.22%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("E38").Value = '  -0.65%  '

$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("E39").Value = '  +13.96%  '

$ws.Range("E40").Value = '  +2.18%  '

$ws.Range("E41").Value = '  -2.86%  '

$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("E42").Value = '  -5.23%  '

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("E43").Value = '  -2.86%  '

$ws.Range("E44").Value = '  -4.87%  '

$ws.Range("E45").Value = '  -2.21%  '

$ws.Range("D46").Value = '2.697.42'
$ws.Range("E46").Value = '  -3.02%  '

$ws.Range("E47").Value = '  -0.73%  '

$ws.Range("E48").Value = '  -5.23%  '

$ws.Range("E49").Value = '  -0.58%  '

$ws.Range("E50").Value = '  -2.63%  '

$ws.Range("E51").Value = '  -1.76%  '

# --- Numeric-looking D-column values forced to remain Text ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.21'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.12'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.599'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.91'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.423'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '30.53'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000174'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.04'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.92'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '382.86'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.86'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.547'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.82'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.79'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000121'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.91'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '24.36'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.90'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.01'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.34'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.19'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.59'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.60'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '29.61'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.80'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.62'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.52'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.47'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0701'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.71'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.74'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0295'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '317.70'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.02'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
